# Clear/replace the STUDYID / DOMAIN / IETESTCD / IETEST / IECAT / IESCAT /
# TIRL / TIVERS header-ish labels that were erroneously duplicated down
# column A/B of the "TI" sheet, and stamp the "TI" domain marker + a
# placeholder space into columns A, B and H for every data row (2-110).
# Row 111 already carries the corrected values and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TI")

for ($r = 2; $r -le 110; $r++) {
    $ws.Range("A$r").Value = " "
    $ws.Range("B$r").Value = "TI"
    $ws.Range("H$r").Value = " "
}
